# Add a new "SheetGradingOrder" worksheet as the first tab of the workbook.
# It lists the grading order (Seq / Sheetname) with a single row pointing
# at the "Sheet3" worksheet, and becomes the active sheet/selection.

$wb = $excel.ActiveWorkbook

# Insert the new sheet before whatever is currently the first tab so it
# ends up as the new first sheet, regardless of which sheet is active.
$firstSheet = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($firstSheet)
$newSheet.Name = "SheetGradingOrder"

$newSheet.Range("A1").Value = "Seq"
$newSheet.Range("B1").Value = "Sheetname"
$newSheet.Range("A2").Value = 1
$newSheet.Range("B2").Value = "Sheet3"

# Make the new sheet active and select A4, matching the saved view state.
$newSheet.Activate() | Out-Null
$newSheet.Range("A4").Select() | Out-Null
